$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "225.37" or
# "33.930.19" are not auto-converted to numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '33.930.19'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.772.26'
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '225.37'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '0.545'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("E7").Value = '  +0.53%  '
$ws.Range("D8").Value = '31.00'
$ws.Range("E8").Value = '  -0.38%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '0.0652'
$ws.Range("E10").Value = '  -1.48%  '
$ws.Range("D11").Value = '0.0927'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '2.029.96'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.800.63'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '10.98'
$ws.Range("E14").Value = '  +9.18%  '
$ws.Range("D15").Value = '0.621'
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("D16").Value = '33.934.88'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '4.18'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").Value = '68.39'
$ws.Range("E18").Value = '  -0.87%  '
$ws.Range("D19").Value = '250.18'
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").Value = '10.26'
$ws.Range("E22").Value = '  -1.38%  '
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("D24").Value = '2.13'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").Value = '155.74'
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").Value = '16.29'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = '0.113'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '6.93'
$ws.Range("E28").Value = '  -1.43%  '
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").Value = '0.0510'
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").Value = '1.82'
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").Value = '1.440.87'
$ws.Range("E35").Value = '  -5.97%  '
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  -1.36%  '
$ws.Range("D37").Value = '0.622'
$ws.Range("E37").Value = '  +1.32%  '
$ws.Range("D38").Value = '0.0185'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("D40").Value = '82.06'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").Value = '0.881'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D43").Value = '2.04'
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("D44").Value = '0.0509'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").Value = '1.927.44'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").Value = '5.76'
$ws.Range("E47").Value = '  +2.19%  '
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("D49").Value = '11.73'
$ws.Range("E49").Value = '  +5.67%  '
$ws.Range("D50").Value = '49.57'
$ws.Range("E50").Value = '  -4.82%  '
$ws.Range("D51").Value = '96.85'
$ws.Range("E51").Value = '  +2.09%  '
